$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("code to write")
$ws.Activate()

# BGR integer colors matching existing conditional fills used on this sheet:
#   "Yes"          -> fill FF92D050 (green)  -> BGR int 5296274
#   "In Progress"  -> fill FFFFFF00 (yellow) -> BGR int 65535
$yesColor = 5296274
$inProgressColor = 65535

# Push the existing rows 27:48 down into 28:49 (row 49 was an unused row
# number in this sheet, so nothing below needs to move), opening up a new
# row 27 for the "ControlParam / XMLCalc" entry.
$ws.Range("A27:C48").Copy() | Out-Null
$ws.Range("A28:C49").PasteSpecial() | Out-Null
$excel.CutCopyMode = $false

# New row: ControlParam / XMLCalc / Yes
$ws.Range("A27").Value = "ControlParam"
$ws.Range("B27").Value = "XMLCalc"
$ws.Range("C27").Value = "Yes"
$ws.Range("C27").Interior.Color = $yesColor

# Update statuses for the existing ControlParam rows
$ws.Range("C25").Value = "In Progress"
$ws.Range("C25").Interior.Color = $inProgressColor

$ws.Range("C26").Value = "Yes"
$ws.Range("C26").Interior.Color = $yesColor

# Rows below shifted down by one after the insert; update their statuses too
$ws.Range("C28").Value = "In Progress"
$ws.Range("C28").Interior.Color = $inProgressColor

$ws.Range("C29").Value = "In Progress"
$ws.Range("C29").Interior.Color = $inProgressColor

$ws.Range("C30").Value = "In Progress"
$ws.Range("C30").Interior.Color = $inProgressColor

$ws.Range("C31").Value = "In Progress"
$ws.Range("C31").Interior.Color = $inProgressColor

# Reflect the author's final view/selection state
$ws.Range("C32").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

# Page setup was touched (portrait / paper size 9 - A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "edit complete"
